$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 413, shifting existing rows 413:442 down to 414:443
$ws.Rows.Item(413).Insert()

# Populate the newly inserted row 413 with the new weekly record
$ws.Cells.Item(413, 1).Value = 4
$ws.Cells.Item(413, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(413, 3).Value = 'Los Lagos'
$ws.Cells.Item(413, 4).Value = 44826
$ws.Cells.Item(413, 5).Value = 10
$ws.Cells.Item(413, 6).Value = 100114013
$ws.Cells.Item(413, 7).Value = 'Zanahoria'
$ws.Cells.Item(413, 8).Value = 'Sin especificar'
$ws.Cells.Item(413, 9).Value = 'Primera'
$ws.Cells.Item(413, 10).Value = 250
$ws.Cells.Item(413, 11).Value = 15000
$ws.Cells.Item(413, 12).Value = 15000
$ws.Cells.Item(413, 13).Value = 15000
$ws.Cells.Item(413, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(413, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(413, 16).Value = 750
$ws.Cells.Item(413, 17).Value = 20
$ws.Cells.Item(413, 18).Value = 'Hortaliza'
